# Handback-status report regeneration: the handoff/handback timestamps for
# the first ("4d66c9be...") and fourth ("e372b22d...") file entries are
# refreshed to reflect a newer report run (they happen to share the same
# timestamp text, as they did before the edit).
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 18:16:15"
$wsZhCn.Range("E5").Value = "2016-03-20 18:16:15"
$wsZhCn.Range("H2").Value = "2016-03-20 18:16:35"
$wsZhCn.Range("H5").Value = "2016-03-20 18:16:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 18:16:18"
$wsDeDe.Range("E5").Value = "2016-03-20 18:16:18"
$wsDeDe.Range("H2").Value = "2016-03-20 18:16:41"
$wsDeDe.Range("H5").Value = "2016-03-20 18:16:41"
